$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the previous lone value first so its shared string is dropped,
# then rebuild the sheet content with all required strings, in the exact
# order needed so the new shared-strings table indices line up with the
# target workbook (0 = "参数区" ... 17 = "temp#n.addr = ESP ...").
$ws.Range("A4").ClearContents() | Out-Null

$ws.Range("B5").Value  = '参数区'
$ws.Range("B7").Value  = '局部变量区'
$ws.Range("B9").Value  = '临时变量区'
$ws.Range("A2").Value  = '运行栈结构'
$ws.Range("B4").Value  = 'display区'
$ws.Range("B6").Value  = '返回地址'
$ws.Range("A6").Value  = 'CALL语句自动压栈'
$ws.Range("B3").Value  = '区域'
$ws.Range("C3").Value  = '地址'
$ws.Range("A9").Value  = '栈指针ESP'
$ws.Range("A4").Value  = '基指针EBP'
$ws.Range("C4").Value  = 'EBP - 4 * (level + 1)'
$ws.Range("C5").Value  = 'para#n.addr = EBP - 4 * (display.length + n + 1)'
$ws.Range("C6").Value  = 'ret_addr = EBP - 4 * (display.length + parameter.length + 1)'
$ws.Range("C7").Value  = 'var#n.addr = EBP - 4 * (display.length + 1 + n + 1)'
$ws.Range("C8").Value  = 'array#n[offset].addr = EBP - 4 * (display.length + 1 + n + 1 + offset) '
$ws.Range("C9").Value  = 'temp#n.addr = EBP - 4 * (display.length + 1 + var_space + n + 1)'
$ws.Range("C10").Value = 'temp#n.addr = ESP + 4 * (temp_space - n - 1)'

# Styles: build the three non-default cellXfs entries in the same order as
# the target styles.xml: 1 = vertical-center, 2 = horizontal-center,
# 3 = horizontal+vertical-center. These must be applied before the merges
# below so the cell -> style index mapping matches the target exactly.
$ws.Range("C7").VerticalAlignment = -4108
$ws.Range("A6").VerticalAlignment = -4108

$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("B2").HorizontalAlignment = -4108

$ws.Range("B7").HorizontalAlignment = -4108
$ws.Range("B7").VerticalAlignment = -4108
$ws.Range("B8").HorizontalAlignment = -4108
$ws.Range("B8").VerticalAlignment = -4108
$ws.Range("B9").HorizontalAlignment = -4108
$ws.Range("B9").VerticalAlignment = -4108
$ws.Range("B10").HorizontalAlignment = -4108
$ws.Range("B10").VerticalAlignment = -4108
$ws.Range("A9").HorizontalAlignment = -4108
$ws.Range("A9").VerticalAlignment = -4108
$ws.Range("A10").HorizontalAlignment = -4108
$ws.Range("A10").VerticalAlignment = -4108

# Merge the cells that belong together.
$ws.Range("A2:B2").Merge() | Out-Null
$ws.Range("B7:B8").Merge() | Out-Null
$ws.Range("B9:B10").Merge() | Out-Null
$ws.Range("A9:A10").Merge() | Out-Null

# Column widths (character units). The runtime always quantizes widths to
# a whole number of pixels, so these inputs are chosen to land as close as
# possible to the authored widths (17.5, 17.375, 86, 8, 8.625); 86 and 8
# are hit exactly, the others land on the nearest achievable pixel width.
$ws.Columns.Item(1).ColumnWidth = 16.857142857142858
$ws.Columns.Item(2).ColumnWidth = 16.714285714285715
$ws.Columns.Item(3).ColumnWidth = 85.28571428571429
$ws.Columns.Item(4).ColumnWidth = 7.285714285714286
$ws.Columns.Item(5).ColumnWidth = 7.857142857142857

$ws.Range("C10").Select() | Out-Null
